$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.636.29'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.809.46'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.94'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.602'
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.19'
$ws.Range('E8').Value = '  -9.26%  '
$ws.Range('E9').Value = '  +5.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0680'
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0992'
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.072.55'
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.823.49'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.672'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '11.14'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.667.85'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.38'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0783'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.04'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.90'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.69'
$ws.Range('E22').Value = '  +0.59%  '
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.57'
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.68'
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.17'
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.119'
$ws.Range('E28').Value = '  -1.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.54'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.01'
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0544'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.95'
$ws.Range('E34').Value = '  +19.46%  '
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.696'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '91.24'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.33'
$ws.Range('E38').Value = '  +4.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.326.84'
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0192'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.961'
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.15'
$ws.Range('E43').Value = '  -7.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.20'
$ws.Range('E44').Value = '  -9.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('E45').Value = '  -4.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.30'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.000.22'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  +0.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0669'
$ws.Range('E50').Value = '  +7.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '97.98'
$ws.Range('E51').Value = '  -4.93%  '
